$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells (row 1): add I0 and IF columns, matching the style of
# the existing header cells (bold, centered, bordered) ---
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data values for columns I (I0) and J (IF), rows 2-80 ---
$I = @(9,7,8,8,9,8,9,8,8,8,7,8,7,7,9,6,7,8,8,8,8,6,9,7,8,8,9,8,7,8,9,8,8,8,8,9,8,8,8,8,7,7,7,7,8,6,9,7,7,8,8,8,8,8,9,8,7,9,7,7,7,7,7,7,7,6,7,9,6,6,7,7,8,4,3,2,6,6,4)
$J = @(9,7,9,8,9,8,9,8,8,8,7,8,7,7,9,7,7,8,8,8,8,6,9,7,8,8,9,8,7,8,9,8,8,8,8,9,8,8,8,8,7,7,7,8,8,6,9,7,7,8,8,8,8,8,9,8,7,9,7,7,7,8,7,7,7,6,7,10,6,6,7,7,8,4,3,2,6,6,4)

for ($i = 0; $i -lt $I.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $I[$i]
    $ws.Cells.Item($row, 10).Value = $J[$i]
}
